$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B1 = 1, B2 = 2 (plain numeric values feeding the Fibonacci-style chain)
$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 2

# B3:B10 share one formula "=B1+B2" (authored with an embedded newline
# before the '+'); Excel stores it as a shared formula (t="shared").
$ws.Range("B3:B10").Formula = "=B1`n+B2"

# The newline in the formula text made Excel auto-grow row 3's height;
# restore the default row height so the row comes back to normal.
$ws.Rows.Item(3).AutoFit()

# Match the saved selection/active cell from the authored workbook.
$ws.Range("B3:B10").Select() | Out-Null
